$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 - Dining Out
$ws.Range("B4").Value = -181022.72
$ws.Range("F4").Value = -181022.72

# Row 5 - Freelance Income
$ws.Range("D5").Value = 338613.39
$ws.Range("F5").Value = 338613.39

# Row 6 - Groceries
$ws.Range("C6").Value = -179006
$ws.Range("F6").Value = -179006

# Row 7 - Interest Income
$ws.Range("D7").Value = 2254401.67
$ws.Range("F7").Value = 2254401.67

# Row 8 - Pets
$ws.Range("B8").Value = -80719.32000000001
$ws.Range("F8").Value = -80719.32000000001

# Row 9 - Pharmacy
$ws.Range("B9").Value = -190509.04
$ws.Range("F9").Value = -190509.04

# Row 10 - Rent
$ws.Range("B10").Value = -80232.52
$ws.Range("F10").Value = -80232.52

# Row 11 - Shopping
$ws.Range("B11").Value = -279746.85
$ws.Range("F11").Value = -279746.85

# Row 12 - Taxes
$ws.Range("B12").Value = -141099.22
$ws.Range("F12").Value = -141099.22

# Row 13 - Transfer From
$ws.Range("E13").Value = 965763.85
$ws.Range("F13").Value = 965763.85

# Row 14 - Transfer To
$ws.Range("E14").Value = -965763.85
$ws.Range("F14").Value = -965763.85

# Row 15 - Utilities
$ws.Range("C15").Value = -83002.44
$ws.Range("F15").Value = -83002.44

# Row 16 - Wages & Salary
$ws.Range("D16").Value = 413324.54
$ws.Range("F16").Value = 413324.54

# Row 17 - Total
$ws.Range("B17").Value = -953329.67
$ws.Range("C17").Value = -262008.44
$ws.Range("D17").Value = 3006339.6
$ws.Range("F17").Value = 1791001.49
